# FORM-11_nodes_config.xlsx - "Updated with LLP Feedbacks"
# Inserts a new "financial_year" field row for the "body_corporates" section
# (mirrors the existing "financial_year" row already present for
# "individual_partners" at row 22), shifting the existing rows 33-51 down
# to 34-52, and updates the sheet view's scroll/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new blank row before row 33, shifting rows 33:51 down to 34:52.
$ws.Rows("33:33").Insert(-4121)  # -4121 = xlShiftDown

# 2) Copy the cell formatting for the new row from the closest matching
#    existing template cells, so the new row reuses the existing style
#    indices instead of Excel cloning brand-new ones.
$ws.Range("A22").Copy()
$ws.Range("A33").PasteSpecial(-4122)  # -4122 = xlPasteFormats

$ws.Range("B22").Copy()
$ws.Range("B33").PasteSpecial(-4122)

$ws.Range("C22").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$ws.Range("D22").Copy()
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("E22").Copy()
$ws.Range("E33").PasteSpecial(-4122)

$ws.Range("F22").Copy()
$ws.Range("F33").PasteSpecial(-4122)

$ws.Range("G32").Copy()
$ws.Range("G33").PasteSpecial(-4122)

$ws.Range("H22").Copy()
$ws.Range("H33").PasteSpecial(-4122)

$ws.Range("I22").Copy()
$ws.Range("I33").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Populate the new row's values.
$ws.Range("A33").Value2 = "financial_year"
$ws.Range("C33").Value2 = "Single"
$ws.Range("D33").Value2 = "ZMCA_LLP_FORM11"
$ws.Range("E33").Value2 = "ANNUAL_RETURN"
$ws.Range("G33").Value2 = "body_corporates"
$ws.Range("H33").Value2 = "financial_year"

# 4) Make sure the row height / formatting flags match the surrounding rows.
$ws.Rows("33:33").RowHeight = $ws.Rows("34:34").RowHeight

# 5) Update the sheet view scroll position / active selection.
$ws.Application.Goto($ws.Range("A7"))
$ws.Range("J28").Select()
